$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the existing Devolucion record's "Archivo" column (B6) that was left blank
$ws.Range("B6").Value = "ServiceDevolucion"

# Add the new ServiceLogIn (login/authentication) records
$ws.Range("A7").Value = "SL1003231009"
$ws.Range("A8").Value = "SL1003231015"
$ws.Range("B7").Value = "ServiceLogIn"
$ws.Range("A9").Value = "SL1003231035"
$ws.Range("A10").Value = "SL1003231113"
$ws.Range("B8").Value = "ServiceLogIn"
$ws.Range("B9").Value = "ServiceLogIn"

# Update the active selection to match the author's saved state
$ws.Range("F12").Select()
